$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6279773116111755
$ws.Range("B1").Value = 0.9818241596221924
$ws.Range("C1").Value = 2.544865846633911
$ws.Range("D1").Value = 6.346468448638916
$ws.Range("E1").Value = 2.121020555496216
